# Refresh the "cryptos" price/volume snapshot (Thu Apr 25 11:34:56 UTC 2024
# GitHub Actions run). Updates Price (D) / Volume-1h (E) text for most rows,
# and swaps the Hedera / NEARProtocol rows (30-31) to match the new ranking.
#
# NOTE: Price values that look like a plain decimal number (e.g. "608.14")
# are prefixed with a leading apostrophe, same as typing them into Excel, so
# they are stored as literal text instead of being auto-converted to a
# number (values such as "63.895.82" already contain two dots so Excel
# can't parse them as numbers and need no such prefix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.895.82'
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").Value = '3.105.31'
$ws.Range("E3").Value = '  -5.36%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '''608.14'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").Value = '''144.82'
$ws.Range("E6").Value = '  -8.78%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '3.101.66'
$ws.Range("E8").Value = '  -5.45%  '
$ws.Range("E9").Value = '  -4.90%  '
$ws.Range("E10").Value = '  -8.09%  '
$ws.Range("D11").Value = '''5.22'
$ws.Range("E11").Value = '  -10.19%  '
$ws.Range("E12").Value = '  -6.17%  '
$ws.Range("D13").Value = '''0.0000248'
$ws.Range("E13").Value = '  -9.58%  '
$ws.Range("E14").Value = '  -10.33%  '
$ws.Range("D15").Value = '3.614.59'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").Value = '63.877.19'
$ws.Range("E17").Value = '  -3.99%  '
$ws.Range("D18").Value = '3.101.23'
$ws.Range("E18").Value = '  -5.26%  '
$ws.Range("E19").Value = '  -8.64%  '
$ws.Range("D20").Value = '''475.83'
$ws.Range("E20").Value = '  -6.24%  '
$ws.Range("D21").Value = '''14.60'
$ws.Range("E21").Value = '  -6.39%  '
$ws.Range("E22").Value = '  -8.05%  '
$ws.Range("D23").Value = '''7.68'
$ws.Range("E23").Value = '  -6.23%  '
$ws.Range("D24").Value = '''13.54'
$ws.Range("E24").Value = '  -8.08%  '
$ws.Range("D25").Value = '''83.33'
$ws.Range("E25").Value = '  -4.26%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -9.60%  '
$ws.Range("E28").Value = '  -9.89%  '
$ws.Range("E29").Value = '  -11.89%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.113'
$ws.Range("E30").Value = '  -13.31%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''6.67'
$ws.Range("E31").Value = '  -5.46%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("D33").Value = '''2.72'
$ws.Range("E33").Value = '  -5.98%  '
$ws.Range("D34").Value = '''26.13'
$ws.Range("E34").Value = '  -7.33%  '
$ws.Range("E35").Value = '  -4.20%  '
$ws.Range("E36").Value = '  -9.37%  '
$ws.Range("D37").Value = '''52.80'
$ws.Range("E37").Value = '  -5.31%  '
$ws.Range("E38").Value = '  -7.41%  '
$ws.Range("D39").Value = '''460.70'
$ws.Range("E39").Value = '  -7.05%  '
$ws.Range("D40").Value = '''2.95'
$ws.Range("E40").Value = '  -14.43%  '
$ws.Range("E41").Value = '  -8.42%  '
$ws.Range("E42").Value = '  -8.56%  '
$ws.Range("E43").Value = '  -5.85%  '
$ws.Range("D44").Value = '2.834.36'
$ws.Range("E44").Value = '  -6.39%  '
$ws.Range("E45").Value = '  -10.13%  '
$ws.Range("E46").Value = '  -12.71%  '
$ws.Range("D47").Value = '''2.41'
$ws.Range("E47").Value = '  -3.64%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '''26.09'
$ws.Range("E49").Value = '  -10.53%  '
$ws.Range("E50").Value = '  -5.71%  '
$ws.Range("D51").Value = '''118.40'
